$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: CLASE, DIA, HORA_INICIO
# Order matters for shared-string table layout: CLASE, then HORA_INICIO, then DIA
$ws.Range("A2").Value = "FIT BOX"
$ws.Range("C2").Value = "20:00"
$ws.Range("B2").Value = "J"

# Update the active selection shown in the sheet view
$ws.Range("B4").Select()
